$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report pair of rows ("Primera"/"Segunda" quality) is
# inserted before the existing row 203, shifting every subsequent row down
# by two. The two new rows duplicate the structure of the other weekly
# entries, using the latest report date (serial 44719).

$ws.Rows.Item(203).Insert()
$ws.Rows.Item(203).Insert()

# Row 203 - "Primera" quality
$ws.Cells.Item(203, 1).Value = 11
$ws.Cells.Item(203, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(203, 3).Value = "Bíobío"
$ws.Cells.Item(203, 4).Value = 44719
$ws.Cells.Item(203, 5).Value = 8
$ws.Cells.Item(203, 6).Value = 100114014
$ws.Cells.Item(203, 7).Value = "Betarraga"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 600
$ws.Cells.Item(203, 11).Value = 600
$ws.Cells.Item(203, 12).Value = 700
$ws.Cells.Item(203, 13).Value = 650
$ws.Cells.Item(203, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(203, 15).Value = "Región Metropolitana"
$ws.Cells.Item(203, 16).Value = 130
$ws.Cells.Item(203, 17).Value = 5
$ws.Cells.Item(203, 18).Value = "Hortaliza"

# Row 204 - "Segunda" quality
$ws.Cells.Item(204, 1).Value = 11
$ws.Cells.Item(204, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(204, 3).Value = "Bíobío"
$ws.Cells.Item(204, 4).Value = 44719
$ws.Cells.Item(204, 5).Value = 8
$ws.Cells.Item(204, 6).Value = 100114014
$ws.Cells.Item(204, 7).Value = "Betarraga"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Segunda"
$ws.Cells.Item(204, 10).Value = 300
$ws.Cells.Item(204, 11).Value = 500
$ws.Cells.Item(204, 12).Value = 500
$ws.Cells.Item(204, 13).Value = 500
$ws.Cells.Item(204, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(204, 15).Value = "Región Metropolitana"
$ws.Cells.Item(204, 16).Value = 100
$ws.Cells.Item(204, 17).Value = 5
$ws.Cells.Item(204, 18).Value = "Hortaliza"

Write-Output "done"
